# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45206 (2023-10-07) to 45208 (2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$ws.Range("C2:C$lastRow").Value = 45208
